# Slide 1, Shape 1 is the title placeholder (ctrTitle) whose text is being
# changed from "定制类对运算符的支持" to "类与运算符", split across three
# runs so that the middle run ("与") is tagged as Chinese (zh-CN) while the
# surrounding runs ("类" and "运算符") stay English (en-GB), matching the
# original run's language.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# Start with the trailing run - it keeps the original run's formatting
# (en-GB) because it's the only run in the text frame at this point.
$tr.Text = "运算符"

# Prepend "与" as its own run (inherits the current formatting, en-GB).
$middle = $tr.InsertBefore("与")

# With exactly two runs ("与","运算符"), TextRange.LanguageID always edits
# the first run in the shape - which is "与" right now - so this correctly
# marks just that run as Chinese.
$sh.TextFrame.TextRange.LanguageID = "zh-CN"

# Prepend "类" as the new first run; it inherits the formatting of what is
# currently the first run ("与", now zh-CN).
$first = $sh.TextFrame.TextRange.InsertBefore("类")

# Re-apply LanguageID: the first run is now "类", so this flips it back to
# en-GB without touching "与" (run 2) or "运算符" (run 3).
$sh.TextFrame.TextRange.LanguageID = "en-GB"
